$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Conversion"
$ws.Range("C2").Value = "Binary"
$ws.Range("A3").Value = 100
$ws.Range("C4").Value = "Binary"
$ws.Range("F8").Value = "Binary"
$ws.Range("A9").Value = "-"
$ws.Range("C9").Value = "'080d"
$ws.Range("C9").NumberFormat = "0.00E+00"
$ws.Range("B9").Value = "081a"
$ws.Range("A5").Value = "0d"

$ws.Range("N6").Select()
